$d = $word.ActiveDocument

# Update stack trace line numbers to reflect Apache POI 4.1.0 -> 5.2.3 move
$d.Content.Find.Execute("ThreadPoolExecutor.java:1130", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ThreadPoolExecutor.java:1136", 2)

$d.Content.Find.Execute("ThreadPoolExecutor`$Worker.run(ThreadPoolExecutor.java:630)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ThreadPoolExecutor`$Worker.run(ThreadPoolExecutor.java:635)", 2)

$d.Content.Find.Execute("Thread.java:832", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Thread.java:833", 2)
